$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply text-number format to all D (Price) and E (Volume) cells in the data rows
# so that numeric-looking strings are preserved exactly as text (matching the
# original inline-string cell content) rather than being coerced into Excel numbers.
$ws.Range("D2:E51").NumberFormat = "@"

# Cell value updates (from the authoritative diff)
$ws.Range('D2').Value2 = '29.300.98'
$ws.Range('E2').Value2 = '  -0.44%  '
$ws.Range('D3').Value2 = '1.840.90'
$ws.Range('E3').Value2 = '  -0.47%  '
$ws.Range('D4').Value2 = '0.9993'
$ws.Range('E4').Value2 = '  -0.08%  '
$ws.Range('D5').Value2 = '240.21'
$ws.Range('E5').Value2 = '  -0.14%  '
$ws.Range('D6').Value2 = '0.6269'
$ws.Range('E6').Value2 = '  -0.13%  '
$ws.Range('D7').Value2 = '0.9999'
$ws.Range('E7').Value2 = '  -0.07%  '
$ws.Range('D8').Value2 = '0.07482'
$ws.Range('E8').Value2 = '  -2.61%  '
$ws.Range('D9').Value2 = '0.2894'
$ws.Range('E9').Value2 = '  -0.93%  '
$ws.Range('D10').Value2 = '24.32'
$ws.Range('E10').Value2 = '  -2.34%  '
$ws.Range('D11').Value2 = '0.07714'
$ws.Range('E11').Value2 = '  -0.46%  '
$ws.Range('D12').Value2 = '1.841.87'
$ws.Range('E12').Value2 = '  -1.13%  '
$ws.Range('D13').Value2 = '4.985'
$ws.Range('E13').Value2 = '  -1.05%  '
$ws.Range('D14').Value2 = '0.6776'
$ws.Range('E14').Value2 = '  -0.61%  '
$ws.Range('D15').Value2 = '0.00001026'
$ws.Range('E15').Value2 = '  -4.78%  '
$ws.Range('D16').Value2 = '82.10'
$ws.Range('E16').Value2 = '  -1.71%  '
$ws.Range('D17').Value2 = '2.101.96'
$ws.Range('E17').Value2 = '  -0.34%  '
$ws.Range('D18').Value2 = '6.102'
$ws.Range('E18').Value2 = '  -1.86%  '
$ws.Range('D19').Value2 = '29.335.92'
$ws.Range('E19').Value2 = '  -0.42%  '
$ws.Range('D20').Value2 = '228.47'
$ws.Range('E20').Value2 = '  -0.15%  '
$ws.Range('D22').Value2 = '1.000'
$ws.Range('E22').Value2 = '  -0.01%  '
$ws.Range('E23').Value2 = '  -1.07%  '
$ws.Range('D24').Value2 = '1.001'
$ws.Range('E24').Value2 = '  -0.05%  '
$ws.Range('D25').Value2 = '158.63'
$ws.Range('E25').Value2 = '  +0.61%  '
$ws.Range('E26').Value2 = '  +0.10%  '
$ws.Range('D27').Value2 = '8.371'
$ws.Range('E27').Value2 = '  -0.56%  '
$ws.Range('D28').Value2 = '17.53'
$ws.Range('E28').Value2 = '  -1.14%  '
$ws.Range('D29').Value2 = '1.392'
$ws.Range('E29').Value2 = '  +1.93%  '
$ws.Range('D30').Value2 = '1.473'
$ws.Range('E30').Value2 = '  +0.77%  '
$ws.Range('E31').Value2 = '  +1.21%  '
$ws.Range('D32').Value2 = '4.094'
$ws.Range('E32').Value2 = '  -0.82%  '
$ws.Range('D33').Value2 = '4.030'
$ws.Range('E33').Value2 = '  -0.52%  '
$ws.Range('D34').Value2 = '1.816'
$ws.Range('E34').Value2 = '  -1.56%  '
$ws.Range('D35').Value2 = '1.142'
$ws.Range('E35').Value2 = '  -1.84%  '
$ws.Range('D36').Value2 = '0.6906'
$ws.Range('E36').Value2 = '  -2.28%  '
$ws.Range('D37').Value2 = '2.585'
$ws.Range('E37').Value2 = '  -0.39%  '
$ws.Range('D38').Value2 = '2.834'
$ws.Range('E38').Value2 = '  +2.73%  '
$ws.Range('D39').Value2 = '1.245.25'
$ws.Range('E39').Value2 = '  +1.64%  '
$ws.Range('D40').Value2 = '0.01813'
$ws.Range('E40').Value2 = '  +1.05%  '
$ws.Range('D41').Value2 = '6.512'
$ws.Range('E41').Value2 = '  +0.82%  '
$ws.Range('D42').Value2 = '0.9056'
$ws.Range('E42').Value2 = '  +0.28%  '
$ws.Range('D43').Value2 = '0.9987'
$ws.Range('E43').Value2 = '  -0.20%  '
$ws.Range('D44').Value2 = '2.003.06'
$ws.Range('E44').Value2 = '  -0.64%  '
$ws.Range('D45').Value2 = '101.22'
$ws.Range('D46').Value2 = '65.68'
$ws.Range('E46').Value2 = '  -0.71%  '
$ws.Range('B47').Value2 = 'Aptos'
$ws.Range('C47').Value2 = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D47').Value2 = '7.074'
$ws.Range('E47').Value2 = '  -1.48%  '
$ws.Range('B48').Value2 = 'Algorand'
$ws.Range('C48').Value2 = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D48').Value2 = '0.1162'
$ws.Range('E48').Value2 = '  +0.37%  '
$ws.Range('B49').Value2 = 'BabyDogeCoin'
$ws.Range('C49').Value2 = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').Value2 = '0.00000000116'
$ws.Range('E49').Value2 = '  -2.67%  '
$ws.Range('D50').Value2 = '9.007'
$ws.Range('E50').Value2 = '  -0.12%  '
$ws.Range('D51').Value2 = '0.3931'
$ws.Range('E51').Value2 = '  -2.25%  '
